$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Introduce the two brand-new table names in the same order the original
# author must have typed them in (movie_link before complete_cast), so the
# shared-string table ends up with matching indices (4=movie_link, 5=complete_cast).
$ws.Range("B4").Value = "movie_link"
$ws.Range("B3").Value = "complete_cast"

# --- Row 3 changes from old (title /movie_info_idx pairing moves away) to new pairing ---
$ws.Range("A3").Value = "title"
$ws.Range("C3").Value = 1000
$ws.Range("D3").Value = 50
$ws.Range("E3").Value = 0.4

# --- New row 4 ---
$ws.Range("A4").Value = "title"
$ws.Range("C4").Value = 1000
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = 0.01

# --- New row 5 (original row 3 content, shifted down) ---
$ws.Range("A5").Value = "movie_info"
$ws.Range("B5").Value = "movie_info_idx"
$ws.Range("C5").Value = 300
$ws.Range("D5").Value = 150
$ws.Range("E5").Value = 0.1

# --- New row 6 ---
$ws.Range("A6").Value = "movie_info"
$ws.Range("B6").Value = "complete_cast"
$ws.Range("C6").Value = 300
$ws.Range("D6").Value = 50
$ws.Range("E6").Value = 0.2

# --- New row 7 ---
$ws.Range("A7").Value = "movie_info"
$ws.Range("B7").Value = "movie_link"
$ws.Range("C7").Value = 300
$ws.Range("D7").Value = 100
$ws.Range("E7").Value = 0.4

# --- New row 8 ---
$ws.Range("A8").Value = "movie_info_idx"
$ws.Range("B8").Value = "complete_cast"
$ws.Range("C8").Value = 150
$ws.Range("D8").Value = 50
$ws.Range("E8").Value = 0.1

# --- New row 9 ---
$ws.Range("A9").Value = "movie_info_idx"
$ws.Range("B9").Value = "movie_link"
$ws.Range("C9").Value = 150
$ws.Range("D9").Value = 100
$ws.Range("E9").Value = 0.05

# --- New row 10 ---
$ws.Range("A10").Value = "complete_cast"
$ws.Range("B10").Value = "movie_link"
$ws.Range("C10").Value = 50
$ws.Range("D10").Value = 100
$ws.Range("E10").Value = 0.2

# --- Column widths (closest achievable to author's 17.7109375 / 17.5703125 OOXML widths) ---
$ws.Columns.Item(1).ColumnWidth = 17
$ws.Columns.Item(2).ColumnWidth = 16.86

# --- View: zoom + selection ---
$excel.ActiveWindow.Zoom = 136
$ws.Range("F10").Select()
